$d = $word.ActiveDocument

# Set the comment author/initials so new comments are attributed correctly.
$word.UserName = "aditya gupta"
$word.UserInitials = "ag"

# ---------------------------------------------------------------------------
# Comment 0: "This needs to be removed"
# Spans from the very start of the document ("Use precise terminology...")
# through the end of the "</think>" paragraph (i.e. across the paragraph
# boundary between paragraph 1 and paragraph 2).
#
# The engine's Comments.Add mis-places the comment markers whenever the
# supplied range crosses a paragraph boundary (it collapses both markers to
# the very start of the document). To work around this we temporarily merge
# paragraph 1 and paragraph 2 into a single paragraph (by deleting the
# paragraph mark between them), add the comment while it's one paragraph
# (so the range no longer crosses a boundary), and then re-split the
# paragraph back apart by re-inserting the paragraph mark at the same spot.
# ---------------------------------------------------------------------------
$pmark = $d.Range(98, 99)
$pmark.Delete()

$c0Range = $d.Range(0, 107)
$c0 = $d.Comments.Add($c0Range, "This needs to be removed")
$c0.Initial = "ag"

$splitPoint = $d.Range(98, 98)
$splitPoint.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Comment 1: "Good but could have added some information about ARAMCO"
# Covers the transaction-narrative paragraph, up to (but not including) its
# very last trailing space character, which is left as its own run after
# the comment reference (matching the target edit).
# ---------------------------------------------------------------------------
$c1Range = $d.Range(435, 1028)
$c1 = $d.Comments.Add($c1Range, "Good but could have added some information about ARAMCO")
$c1.Initial = "ag"

# ---------------------------------------------------------------------------
# Comment 2: "Good covers expected activities"
# Wraps the entire "JDF Industries (Customer ID: C-4) ..." paragraph.
# ---------------------------------------------------------------------------
$c2Range = $d.Range(1031, 1548)
$c2 = $d.Comments.Add($c2Range, "Good covers expected activities")
$c2.Initial = "ag"

# ---------------------------------------------------------------------------
# Comment 3: "This is a great point"
# Wraps the entire "4. Mirroring transaction patterns ..." paragraph.
# ---------------------------------------------------------------------------
$c3Range = $d.Range(2038, 2149)
$c3 = $d.Comments.Add($c3Range, "This is a great point")
$c3.Initial = "ag"

Write-Output "Comments added: $($d.Comments.Count)"
